$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update matched-error values for the ifoCAST full series evaluation (values shifted
# by one row since a new realized observation now factors into every horizons error).
$ws.Range("B2").Value = 0.8032788223262832
$ws.Range("C2").Value = -1.002366475608788
$ws.Range("D2").Value = -1.464557645652444
$ws.Range("E2").Value = 0.9087989414264609
$ws.Range("F2").Value = 0.03490493049987109
$ws.Range("G2").Value = 0.2751479164465338
$ws.Range("H2").Value = -0.04060839016473031
$ws.Range("I2").Value = 0.4714428511361884
$ws.Range("J2").Value = 0.1848835897783218
$ws.Range("K2").Value = 0.2125075656625323
$ws.Range("B3").Value = -0.5941340219870179
$ws.Range("C3").Value = -1.21293460176605
$ws.Range("D3").Value = 0.8268928876033257
$ws.Range("E3").Value = 0.1099969798811567
$ws.Range("F3").Value = 0.3810550769602576
$ws.Range("G3").Value = 0.0178082559939251
$ws.Range("H3").Value = 0.5378815349662799
$ws.Range("I3").Value = 0.2595655483891583
$ws.Range("J3").Value = 0.2828208575635111
$ws.Range("K3").Value = 0.3749895042266514
$ws.Range("B4").Value = -1.318433813614865
$ws.Range("C4").Value = 0.7007876036379678
$ws.Range("D4").Value = 0.2782450485143884
$ws.Range("E4").Value = 0.4417191695641399
$ws.Range("F4").Value = 0.02894109536855799
$ws.Range("G4").Value = 0.6007442862932105
$ws.Range("H4").Value = 0.3169170761829015
$ws.Range("I4").Value = 0.3273260345678901
$ws.Range("J4").Value = 0.42606337525307
$ws.Range("K4").Value = 0.6412619431822899
$ws.Range("B5").Value = 0.0963469837902291
$ws.Range("C5").Value = 0.214041671159695
$ws.Range("D5").Value = 0.7332664776213567
$ws.Range("E5").Value = -0.008765780717375604
$ws.Range("F5").Value = 0.5921602662197494
$ws.Range("G5").Value = 0.4006581537802698
$ws.Range("H5").Value = 0.3621531794959351
$ws.Range("I5").Value = 0.452627758305367
$ws.Range("J5").Value = 0.6867607839288887
$ws.Range("K5").Value = 0.1781578843816368
$ws.Range("B6").Value = 0.1264008423207837
$ws.Range("C6").Value = 0.6808472755916881
$ws.Range("D6").Value = 0.04988061626763002
$ws.Range("E6").Value = 0.5970894115568507
$ws.Range("F6").Value = 0.3895648707313746
$ws.Range("G6").Value = 0.3749198787210216
$ws.Range("H6").Value = 0.4613573173527261
$ws.Range("I6").Value = 0.6892627280777406
$ws.Range("J6").Value = 0.1844338218533179
$ws.Range("K6").Value = 0.4715052544735016
$ws.Range("B7").Value = 0.7032752552246967
$ws.Range("C7").Value = 0.0524841558300787
$ws.Range("D7").Value = 0.5750217259028355
$ws.Range("E7").Value = 0.3835828904270196
$ws.Range("F7").Value = 0.3695766285386105
$ws.Range("G7").Value = 0.4509887839823598
$ws.Range("H7").Value = 0.6807854059541167
$ws.Range("I7").Value = 0.1766615135465071
$ws.Range("J7").Value = 0.4628812809405329
$ws.Range("K7").Value = 0.2074249537672726
$ws.Range("B8").Value = 0.0739087272872988
$ws.Range("C8").Value = 0.6951995747020479
$ws.Range("D8").Value = 0.2943885460132365
$ws.Range("E8").Value = 0.3398496246900327
$ws.Range("F8").Value = 0.4596149040122699
$ws.Range("G8").Value = 0.6548452325286815
$ws.Range("H8").Value = 0.1538829505182796
$ws.Range("I8").Value = 0.4492632457919151
$ws.Range("J8").Value = 0.1889041565820968
$ws.Range("K8").Value = 0.5706702220727796
$ws.Range("B9").Value = 0.6528789423816584
$ws.Range("C9").Value = 0.2651840721575033
$ws.Range("D9").Value = 0.3368290248851115
$ws.Range("E9").Value = 0.4407536204007895
$ws.Range("F9").Value = 0.6345141014634773
$ws.Range("G9").Value = 0.1393527950840318
$ws.Range("H9").Value = 0.4326222002996472
$ws.Range("I9").Value = 0.1713203111533466
$ws.Range("J9").Value = 0.5541963385427369
$ws.Range("K9").Value = 0.1437698493309027
$ws.Range("B10").Value = 0.6053818127754134
$ws.Range("C10").Value = 0.4122000866690486
$ws.Range("D10").Value = 0.2496603340877904
$ws.Range("E10").Value = 0.6592080140502106
$ws.Range("F10").Value = 0.1487321986403278
$ws.Range("G10").Value = 0.3778114016882561
$ws.Range("H10").Value = 0.1524262202646768
$ws.Range("I10").Value = 0.5393323377276911
$ws.Range("J10").Value = 0.115058138701532
$ws.Range("K10").Value = 0.4067718394308724
$ws.Range("B11").Value = 0.8628949586592991
$ws.Range("C11").Value = 0.2967710363001488
$ws.Range("D11").Value = 0.4189247832594023
$ws.Range("E11").Value = 0.1846772797061906
$ws.Range("F11").Value = 0.3757606442486632
$ws.Range("G11").Value = 0.07585798082864662
$ws.Range("H11").Value = 0.5119329433524077
$ws.Range("I11").Value = 0.08876908850380663
$ws.Range("J11").Value = 0.3633745487175398
$ws.Range("B12").Value = 0.6090966232236873
$ws.Range("C12").Value = 0.5522135229949265
$ws.Range("D12").Value = 0.005598857889999004
$ws.Range("E12").Value = 0.4039548830192304
$ws.Range("F12").Value = 0.1118832920210401
$ws.Range("G12").Value = 0.4743913731481941
$ws.Range("H12").Value = 0.08322674941644539
$ws.Range("I12").Value = 0.3675498776562884
$ws.Range("B13").Value = 0.7878040141027678
$ws.Range("C13").Value = 0.09027759876430858
$ws.Range("D13").Value = 0.2583545163855133
$ws.Range("E13").Value = 0.1218370348802827
$ws.Range("F13").Value = 0.489756542847739
$ws.Range("G13").Value = 0.04437841445902233
$ws.Range("H13").Value = 0.3517040686291025
$ws.Range("B14").Value = 0.4013017852456914
$ws.Range("C14").Value = 0.3754432907967085
$ws.Range("D14").Value = -0.04062710656928412
$ws.Range("E14").Value = 0.5187154933129405
$ws.Range("F14").Value = 0.08012128691392592
$ws.Range("G14").Value = 0.3203764222454754
$ws.Range("B15").Value = 0.6222684682008229
$ws.Range("C15").Value = -0.02297123903139461
$ws.Range("D15").Value = 0.4240932542019461
$ws.Range("E15").Value = 0.112338675162406
$ws.Range("F15").Value = 0.3327645480731927
$ws.Range("B16").Value = 0.2167051203848173
$ws.Range("C16").Value = 0.5091174976711597
$ws.Range("D16").Value = -0.004145903195608092
$ws.Range("E16").Value = 0.3478698197250452
$ws.Range("B17").Value = 0.6739775747052469
$ws.Range("C17").Value = 0.009391369052308113
$ws.Range("D17").Value = 0.2848969007350822
$ws.Range("B18").Value = 0.2632404109177161
$ws.Range("C18").Value = 0.3842149509171186
$ws.Range("B19").Value = 0.4282746421565676

# Clear the trailing cells that fall outside the (now one-row-shorter) staircase of
# available forecast horizons for the most recent observations.
$ws.Range("K11").ClearContents()
$ws.Range("J12").ClearContents()
$ws.Range("I13").ClearContents()
$ws.Range("H14").ClearContents()
$ws.Range("G15").ClearContents()
$ws.Range("F16").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("B20").ClearContents()
